# Update roster dropout dates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Candidates whose "Date Dropped" cell currently holds the placeholder "-"
# string are being given real drop-out dates. Column G is "Date Dropped".
#   Row 2  -> Bennet     -> 2020-02-11
#   Row 8  -> Buttigieg  -> 2020-03-01
#   Row 17 -> Klobuchar  -> 2020-03-02
#   Row 23 -> Yang       -> 2020-02-11

# Copy the date number format already used elsewhere in the column (e.g.
# G3) onto the cells being updated, so they share the same style as the
# other dated rows instead of getting a brand-new style definition.
$ws.Range("G3").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G8").PasteSpecial(-4122)
$ws.Range("G17").PasteSpecial(-4122)
$ws.Range("G23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("G2").Value = [DateTime]"2020-02-11"
$ws.Range("G8").Value = [DateTime]"2020-03-01"
$ws.Range("G17").Value = [DateTime]"2020-03-02"
$ws.Range("G23").Value = [DateTime]"2020-02-11"

# Move the active selection, matching the author's final cursor position.
$ws.Range("I18").Select()
